$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.723.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.888.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.82%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.886.44'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.27%  '
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.427'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000230'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.39%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.368.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.691.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.892.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  -11.03%  '
$ws.Range("E28").Value = '  -5.12%  '
$ws.Range("E29").Value = '  +8.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("E31").Value = '  -4.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.74%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.956'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("E39").Value = '  -4.62%  '
$ws.Range("E40").Value = '  -6.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("E44").Value = '  -4.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.677.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '343.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.81%  '
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.84%  '
